$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/12/2024  Through  2/18/2024"

# --- Cells whose type/style changes (text <-> number) ---
# Strategy: set raw value first (apostrophe-prefixed for numeric-looking text),
# then PasteSpecial(xlPasteFormats=-4122) from a same-style donor cell so the
# resulting style index matches the target (text style 14, integer style 15,
# percent style 16) without fabricating a brand-new style entry.
$ws.Range("D14").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("C15").Value = 2
$ws.Range("G14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("G14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = 100
$ws.Range("H14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F15").Value = 2
$ws.Range("G14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("I15").Value = 2
$ws.Range("G14").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("G14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1
$ws.Range("G14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("C26").Value = 2
$ws.Range("G14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$ws.Range("G14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = 100
$ws.Range("H14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("C27").Value = 2
$ws.Range("G14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D28").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("D29").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Plain value updates (style/type unchanged) ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -60
$ws.Range("M15").Value = -33.333333333333
$ws.Range("N15").Value = -87.5
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 44.444444444444
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 30
$ws.Range("L16").Value = -40.90909090909
$ws.Range("M16").Value = -36.585365853658
$ws.Range("N16").Value = -91.186440677966
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -4
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 46
$ws.Range("K17").Value = -13.043478260869
$ws.Range("L17").Value = -9.090909090909
$ws.Range("M17").Value = -13.043478260869
$ws.Range("N17").Value = -61.904761904761
$ws.Range("D18").Value = 2
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -30
$ws.Range("J18").Value = 22
$ws.Range("K18").Value = -40.90909090909
$ws.Range("L18").Value = -35
$ws.Range("M18").Value = -58.064516129032
$ws.Range("N18").Value = -97.23991507431
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -22.222222222222
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 53
$ws.Range("J19").Value = 88
$ws.Range("K19").Value = -39.772727272727
$ws.Range("L19").Value = -30.263157894736
$ws.Range("M19").Value = -19.696969696969
$ws.Range("N19").Value = -58.91472868217
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 10
$ws.Range("H20").Value = 42.857142857142
$ws.Range("I20").Value = 18
$ws.Range("J20").Value = 15
$ws.Range("K20").Value = 20
$ws.Range("L20").Value = -30.76923076923
$ws.Range("M20").Value = -35.714285714285
$ws.Range("N20").Value = -94.658753709198
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -10.526315789473
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 102
$ws.Range("H21").Value = -9.803921568627
$ws.Range("I21").Value = 152
$ws.Range("J21").Value = 194
$ws.Range("K21").Value = -21.649484536082
$ws.Range("L21").Value = -29.302325581395
$ws.Range("M21").Value = -29.302325581395
$ws.Range("N21").Value = -88.839941262848
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 3
$ws.Range("L22").Value = 50
$ws.Range("M22").Value = -25
$ws.Range("C24").Value = 50
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 108.333333333333
$ws.Range("F24").Value = 157
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = 23.622047244094
$ws.Range("I24").Value = 261
$ws.Range("J24").Value = 207
$ws.Range("K24").Value = 26.086956521739
$ws.Range("L24").Value = 82.517482517482
$ws.Range("M24").Value = 69.480519480519
$ws.Range("C25").Value = 11
$ws.Range("E25").Value = 22.222222222222
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 25
$ws.Range("I25").Value = 78
$ws.Range("J25").Value = 67
$ws.Range("K25").Value = 16.417910447761
$ws.Range("L25").Value = 11.428571428571
$ws.Range("M25").Value = -26.415094339622
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 4
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = -20
$ws.Range("L26").Value = -33.333333333333
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 7
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = 40
$ws.Range("L27").Value = 0
